$d = $word.ActiveDocument

# Locate the "GIS & Geospatial Analysis Consulting" paragraph under the
# Siege Analytics / PARTNER job heading — the new bullet points get
# inserted directly after it.
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "GIS & Geospatial Analysis Consulting") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find 'GIS & Geospatial Analysis Consulting' paragraph"
}

$newBullets = @(
    "• Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels",
    "• Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide",
    "• Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis"
)

$anchor = $d.Paragraphs.Item($targetIndex)
foreach ($bulletText in $newBullets) {
    $anchor.Range.InsertParagraphAfter()
    $targetIndex = $targetIndex + 1
    $newPara = $d.Paragraphs.Item($targetIndex)
    $newPara.Range.Text = $bulletText
    $anchor = $newPara
}

Write-Output "Inserted $($newBullets.Count) new bullet paragraphs after 'GIS & Geospatial Analysis Consulting'."
